$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values (regenerated s_val data, filtering save games)
$data = @{
    2 = @{ B = 0.3464964993005633; C = 9.226618575922256;  D = 3.082599426703578;  E = 6.48142807727062;  G = 19.13714257919702 }
    3 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 16.98373111632243;  E = 6.48142807727062;  G = 28.30127388105354 }
    4 = @{ B = 1.505614041169197;  C = 0.3375848360084654; D = 0.7127328510149897; E = 0.4998867070740569; G = 3.055818435266709 }
    5 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 3.082599426703578;  E = 0.4998867070740569; G = 8.418600821238126 }
    6 = @{ B = 0.7287194209349384; C = 1.65323645889881;   D = 0.1529057820181812; E = 0.4998867070740569; G = 3.034748368925986 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    $ws.Range("B$row").Value = $cols.B
    $ws.Range("C$row").Value = $cols.C
    $ws.Range("D$row").Value = $cols.D
    $ws.Range("E$row").Value = $cols.E
    $ws.Range("G$row").Value = $cols.G
}
